$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115, shifting existing rows 115-242 down to 116-243.
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new weekly record.
$ws.Cells.Item(115, 1).Value = 4
$ws.Cells.Item(115, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(115, 3).Value = "Los Lagos"
$ws.Cells.Item(115, 4).Value = 44539
$ws.Cells.Item(115, 5).Value = 10
$ws.Cells.Item(115, 6).Value = 100112008
$ws.Cells.Item(115, 7).Value = "Coliflor"
$ws.Cells.Item(115, 8).Value = "Sin especificar"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 500
$ws.Cells.Item(115, 11).Value = 1000
$ws.Cells.Item(115, 12).Value = 1200
$ws.Cells.Item(115, 13).Value = 1100
$ws.Cells.Item(115, 14).Value = "$/unidad"
$ws.Cells.Item(115, 15).Value = "Región Metropolitana"
$ws.Cells.Item(115, 16).Value = 1100
$ws.Cells.Item(115, 17).Value = 1
$ws.Cells.Item(115, 18).Value = "Hortaliza"

# Apply the same date-cell number format (style index "2") used by the other rows in column D.
$ws.Cells.Item(115, 4).NumberFormat = $ws.Cells.Item(116, 4).NumberFormat
